$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H values (week of 4/25/2020), matching column G's date style
$ws.Range("H1").Value = 43946
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$values = @{
    2  = 28045
    3  = 17381
    4  = 211
    5  = 6311
    6  = 14419
    7  = 6236
    8  = 1703
    9  = 24866
    10 = 9771
    11 = 1109
    12 = 10053
    13 = 104
    14 = 10327
    15 = 7913
    16 = 1636
    17 = 3337
    18 = 25049
    19 = 14483
    20 = 7470
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 8).Value = $values[$row]
}

# Match column H width to column G (best available precision)
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Update the selected cell
$ws.Range("E9").Select()

